$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Implement enemies/npcs"
$ws.Range("D15").Value = 6
$ws.Range("D14").Value = 14

$ws.Range("A3").Select()
